# Daily attendance processing - 2025-12-12 20:30:13
# Reverse the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
